# Add two new columns I (I0) and J (IF) to the sheet, mirroring the
# existing header/data layout (columns B..H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the other header cells (bold,
# bordered, centered) by copying the format from H1 onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-26) ---
$data = @{
    2  = @(1, 4)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 6)
    6  = @(1, 4)
    7  = @(4, 8)
    8  = @(8, 8)
    9  = @(3, 5)
    10 = @(3, 9)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(2, 6)
    14 = @(1, 6)
    15 = @(3, 7)
    16 = @(3, 7)
    17 = @(1, 6)
    18 = @(2, 6)
    19 = @(6, 8)
    20 = @(6, 8)
    21 = @(5, 8)
    22 = @(1, 5)
    23 = @(4, 7)
    24 = @(4, 5)
    25 = @(1, 2)
    26 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
